$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2023-09-16 Saturday" "2023-09-17 Sunday"

Replace-Text "27÷3=" "55÷9="
Replace-Text "39÷3=" "54÷5="
Replace-Text "57÷3=" "79÷5="
Replace-Text "39÷9=" "77÷5="
Replace-Text "83÷2=" "77÷6="

Replace-Text "11÷7=" "72÷8="
Replace-Text "41÷4=" "38÷4="
Replace-Text "44÷4=" "71÷2="
Replace-Text "91÷6=" "92÷4="
Replace-Text "45÷3=" "47÷9="

Replace-Text "30÷3=" "20÷3="
Replace-Text "12÷9=" "76÷4="
Replace-Text "79÷7=" "88÷6="
Replace-Text "33÷5=" "93÷6="
Replace-Text "88÷7=" "72÷4="

Replace-Text "27÷6=" "86÷2="
Replace-Text "60÷6=" "86÷5="
Replace-Text "44÷2=" "89÷9="
Replace-Text "46÷8=" "86÷8="
Replace-Text "15÷8=" "28÷5="

Replace-Text "80÷3=" "75÷9="
Replace-Text "68÷9=" "49÷9="
Replace-Text "51÷5=" "74÷3="
Replace-Text "56÷4=" "36÷8="
Replace-Text "31÷6=" "72÷3="
